$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 38535.848
$ws.Range("I6").Value = 80.5
$ws.Range("J6").Value = 500000
$ws.Range("K6").Value = 241.5
$ws.Range("L6").Value = 1500000
$ws.Range("M6").Value = -129.5
$ws.Range("N6").Value = -1500224

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 828.2857
$ws.Range("I97").Value = 400
$ws.Range("J97").Value = 899.6667
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 2699.0001
$ws.Range("M97").Value = -704
$ws.Range("N97").Value = -3691.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2962.96
$ws.Range("I137").Value = 1344.1177
$ws.Range("J137").Value = 6403
$ws.Range("K137").Value = 4032.3531
$ws.Range("L137").Value = 19209
$ws.Range("M137").Value = -1482.3531
$ws.Range("N137").Value = -24309

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3448.7
$ws.Range("I138").Value = 2623.5833
$ws.Range("J138").Value = 3709.2632
$ws.Range("K138").Value = 7870.749899999999
$ws.Range("L138").Value = 11127.7896
$ws.Range("M138").Value = -2730.749899999999
$ws.Range("N138").Value = -21407.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 35715852
$ws.Range("I74").Value = 43479640
$ws.Range("J74").Value = 2428.2
$ws.Range("K74").Value = 43479640
$ws.Range("L74").Value = 2428.2
$ws.Range("M74").Value = -43478766
$ws.Range("N74").Value = -4176.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 35715852
$ws.Range("I77").Value = 43479640
$ws.Range("J77").Value = 2428.2
$ws.Range("K77").Value = 217398200
$ws.Range("L77").Value = 12141
$ws.Range("M77").Value = -217393832
$ws.Range("N77").Value = -20877

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 225207.56
$ws.Range("I102").Value = 335736.5
$ws.Range("J102").Value = 4149.6665
$ws.Range("K102").Value = 335736.5
$ws.Range("L102").Value = 4149.6665
$ws.Range("M102").Value = -334114.5
$ws.Range("N102").Value = -7393.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3441.0286
$ws.Range("I110").Value = 2395.6072
$ws.Range("J110").Value = 7622.7144
$ws.Range("K110").Value = 2395.6072
$ws.Range("L110").Value = 7622.7144
$ws.Range("M110").Value = -350.6071999999999
$ws.Range("N110").Value = -11712.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1720.7646
$ws.Range("I99").Value = 889.44446
$ws.Range("J99").Value = 2656
$ws.Range("K99").Value = 889.44446
$ws.Range("L99").Value = 2656
$ws.Range("M99").Value = 608.55554
$ws.Range("N99").Value = -5652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2440.2942
$ws.Range("I107").Value = 1800.5
$ws.Range("J107").Value = 2789.2727
$ws.Range("K107").Value = 1800.5
$ws.Range("L107").Value = 2789.2727
$ws.Range("M107").Value = 119.5
$ws.Range("N107").Value = -6629.2727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 2999.5
$ws.Range("I128").Value = 2999.5
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 8998.5
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -6508.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2844.8076
$ws.Range("I134").Value = 2237.6316
$ws.Range("J134").Value = 4492.857
$ws.Range("K134").Value = 6712.8948
$ws.Range("L134").Value = 13478.571
$ws.Range("M134").Value = -4177.8948
$ws.Range("N134").Value = -18548.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4904063.5
$ws.Range("I31").Value = 1886.9
$ws.Range("J31").Value = 8066757.5
$ws.Range("K31").Value = 1886.9
$ws.Range("L31").Value = 8066757.5
$ws.Range("M31").Value = -1591.9
$ws.Range("N31").Value = -8067347.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4904063.5
$ws.Range("I34").Value = 1886.9
$ws.Range("J34").Value = 8066757.5
$ws.Range("K34").Value = 1886.9
$ws.Range("L34").Value = 8066757.5
$ws.Range("M34").Value = -1684.9
$ws.Range("N34").Value = -8067161.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 35995
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 35995
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 35995
$ws.Range("N48").Value = -36947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 21591.6
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 32652.666
$ws.Range("K57").Value = 5000
$ws.Range("L57").Value = 32652.666
$ws.Range("M57").Value = -4440
$ws.Range("N57").Value = -33772.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 6261.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 6261.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 6261.5
$ws.Range("N96").Value = -11753.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 51073.562
$ws.Range("I132").Value = 56281.594
$ws.Range("J132").Value = 2899.25
$ws.Range("K132").Value = 168844.782
$ws.Range("L132").Value = 8697.75
$ws.Range("M132").Value = -166314.782
$ws.Range("N132").Value = -13757.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36410844
$ws.Range("I4").Value = 48812980
$ws.Range("J4").Value = 20563666
$ws.Range("K4").Value = 146438940
$ws.Range("L4").Value = 61690998
$ws.Range("M4").Value = -146438828
$ws.Range("N4").Value = -61691222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 329
$ws.Range("I26").Value = 121.666664
$ws.Range("J26").Value = 640
$ws.Range("K26").Value = 364.999992
$ws.Range("L26").Value = 1920
$ws.Range("M26").Value = -76.99999200000002
$ws.Range("N26").Value = -2496

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1149.75
$ws.Range("I92").Value = 900
$ws.Range("J92").Value = 1899
$ws.Range("K92").Value = 2700
$ws.Range("L92").Value = 5697
$ws.Range("M92").Value = -1452
$ws.Range("N92").Value = -8193

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1004.8182
$ws.Range("I97").Value = 786.4
$ws.Range("J97").Value = 1472.8572
$ws.Range("K97").Value = 2359.2
$ws.Range("L97").Value = 4418.571599999999
$ws.Range("M97").Value = -1863.2
$ws.Range("N97").Value = -5410.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3469
$ws.Range("I113").Value = 1408.6
$ws.Range("J113").Value = 4940.7144
$ws.Range("K113").Value = 1408.6
$ws.Range("L113").Value = 4940.7144
$ws.Range("M113").Value = 761.4000000000001
$ws.Range("N113").Value = -9280.714400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 66467.31
$ws.Range("I126").Value = 103347.7
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 310043.1
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -307573.1
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6594.4
$ws.Range("I7").Value = 4539.1
$ws.Range("J7").Value = 8649.700000000001
$ws.Range("K7").Value = 4539.1
$ws.Range("L7").Value = 8649.700000000001
$ws.Range("M7").Value = -4427.1
$ws.Range("N7").Value = -8873.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3363.2727
$ws.Range("I22").Value = 2513.8572
$ws.Range("J22").Value = 4849.75
$ws.Range("K22").Value = 2513.8572
$ws.Range("L22").Value = 4849.75
$ws.Range("M22").Value = -2218.8572
$ws.Range("N22").Value = -5439.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3363.2727
$ws.Range("I27").Value = 2513.8572
$ws.Range("J27").Value = 4849.75
$ws.Range("K27").Value = 2513.8572
$ws.Range("L27").Value = 4849.75
$ws.Range("M27").Value = -2406.8572
$ws.Range("N27").Value = -5063.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8545
$ws.Range("I40").Value = 8730.143
$ws.Range("J40").Value = 8415.4
$ws.Range("K40").Value = 8730.143
$ws.Range("L40").Value = 8415.4
$ws.Range("M40").Value = -8594.143
$ws.Range("N40").Value = -8687.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3461.9092
$ws.Range("I61").Value = 2132.3333
$ws.Range("J61").Value = 5057.4
$ws.Range("K61").Value = 2132.3333
$ws.Range("L61").Value = 5057.4
$ws.Range("M61").Value = -1930.3333
$ws.Range("N61").Value = -5461.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 65641.336
$ws.Range("I74").Value = 24990
$ws.Range("J74").Value = 85967
$ws.Range("K74").Value = 24990
$ws.Range("L74").Value = 85967
$ws.Range("M74").Value = -23992
$ws.Range("N74").Value = -87963

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 65641.336
$ws.Range("I77").Value = 24990
$ws.Range("J77").Value = 85967
$ws.Range("K77").Value = 74970
$ws.Range("L77").Value = 257901
$ws.Range("M77").Value = -69978
$ws.Range("N77").Value = -267885

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 46694.5
$ws.Range("I92").Value = 43000
$ws.Range("J92").Value = 50389
$ws.Range("K92").Value = 43000
$ws.Range("L92").Value = 50389
$ws.Range("M92").Value = -40504
$ws.Range("N92").Value = -55381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 324114.9
$ws.Range("I93").Value = 1159.5834
$ws.Range("J93").Value = 1381059.5
$ws.Range("K93").Value = 1159.5834
$ws.Range("L93").Value = 1381059.5
$ws.Range("M93").Value = 88.41660000000002
$ws.Range("N93").Value = -1383555.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3461.9092
$ws.Range("I113").Value = 2132.3333
$ws.Range("J113").Value = 5057.4
$ws.Range("K113").Value = 2132.3333
$ws.Range("L113").Value = 5057.4
$ws.Range("M113").Value = 37.66670000000022
$ws.Range("N113").Value = -9397.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6594.4
$ws.Range("I126").Value = 4539.1
$ws.Range("J126").Value = 8649.700000000001
$ws.Range("K126").Value = 13617.3
$ws.Range("L126").Value = 25949.1
$ws.Range("M126").Value = -11147.3
$ws.Range("N126").Value = -30889.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1336705.2
$ws.Range("I136").Value = 2002818.1
$ws.Range("J136").Value = 4479.6
$ws.Range("K136").Value = 6008454.300000001
$ws.Range("L136").Value = 13438.8
$ws.Range("M136").Value = -6005904.300000001
$ws.Range("N136").Value = -18538.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 24805.5
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 24805.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 24805.5
$ws.Range("N47").Value = -25949.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12361.728
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 14108.777
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 14108.777
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -15356.777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 12361.728
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 14108.777
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 70543.88499999999
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -76783.88499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 37650
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 37650
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 37650
$ws.Range("N80").Value = -39646

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3744
$ws.Range("I81").Value = 2274.1333
$ws.Range("J81").Value = 6500
$ws.Range("K81").Value = 4548.2666
$ws.Range("L81").Value = 13000
$ws.Range("M81").Value = -3487.2666
$ws.Range("N81").Value = -15122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 10701
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 10701
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 10701
$ws.Range("N82").Value = -11467

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 37650
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 37650
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 112950
$ws.Range("N83").Value = -122934

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3744
$ws.Range("I84").Value = 2274.1333
$ws.Range("J84").Value = 6500
$ws.Range("K84").Value = 22741.333
$ws.Range("L84").Value = 65000
$ws.Range("M84").Value = -17437.333
$ws.Range("N84").Value = -75608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 10701
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 10701
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 10701
$ws.Range("N85").Value = -13353
